$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 297.16666
$ws.Range("I2").Value = 259.2
$ws.Range("J2").Value = 324.2857
$ws.Range("K2").Value = 259.2
$ws.Range("L2").Value = 324.2857
$ws.Range("M2").Value = -146.2
$ws.Range("N2").Value = -550.2857
$ws.Range("H32").Value = 1466.3334
$ws.Range("J32").Value = 949.5
$ws.Range("L32").Value = 949.5
$ws.Range("N32").Value = -1601.5
$ws.Range("H33").Value = 280.6316
$ws.Range("J33").Value = 113.42857
$ws.Range("L33").Value = 113.42857
$ws.Range("N33").Value = -571.42857
$ws.Range("H64").Value = 3218.75
$ws.Range("I64").Value = 2958.3333
$ws.Range("K64").Value = 2958.3333
$ws.Range("M64").Value = -2710.3333
$ws.Range("H67").Value = 3218.75
$ws.Range("I67").Value = 2958.3333
$ws.Range("K67").Value = 2958.3333
$ws.Range("M67").Value = -2100.3333
$ws.Range("H138").Value = 3913.9736
$ws.Range("I138").Value = 7725.857
$ws.Range("J138").Value = 3053.2258
$ws.Range("K138").Value = 23177.571
$ws.Range("L138").Value = 9159.6774
$ws.Range("M138").Value = -18037.571
$ws.Range("N138").Value = -19439.6774
$ws.Range("H141").Value = 1274495.4
$ws.Range("I141").Value = 1556455.5
$ws.Range("K141").Value = 4669366.5
$ws.Range("M141").Value = -4664186.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1339.1562
$ws.Range("I45").Value = 1131.5714
$ws.Range("K45").Value = 1131.5714
$ws.Range("M45").Value = -754.5714
$ws.Range("H74").Value = 1933.7693
$ws.Range("I74").Value = 1783
$ws.Range("K74").Value = 1783
$ws.Range("M74").Value = -909
$ws.Range("H77").Value = 1933.7693
$ws.Range("I77").Value = 1783
$ws.Range("K77").Value = 8915
$ws.Range("M77").Value = -4547
$ws.Range("H122").Value = 1835.3158
$ws.Range("I122").Value = 1826.1666
$ws.Range("K122").Value = 5478.4998
$ws.Range("M122").Value = -3028.4998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H86").Value = 136538.94
$ws.Range("I86").Value = 4726.857
$ws.Range("J86").Value = 251874.5
$ws.Range("K86").Value = 4726.857
$ws.Range("L86").Value = 251874.5
$ws.Range("M86").Value = -3603.857
$ws.Range("N86").Value = -254120.5
$ws.Range("H89").Value = 136538.94
$ws.Range("I89").Value = 4726.857
$ws.Range("J89").Value = 251874.5
$ws.Range("K89").Value = 23634.285
$ws.Range("L89").Value = 1259372.5
$ws.Range("M89").Value = -18018.285
$ws.Range("N89").Value = -1270604.5
$ws.Range("H107").Value = 2002.45
$ws.Range("I107").Value = 1778.4
$ws.Range("K107").Value = 1778.4
$ws.Range("M107").Value = 141.5999999999999
$ws.Range("N30").ClearContents()  # cell removed entirely in target (no M30/M140-style neighbor either)

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1500
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1666.6666
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1666.6666
$ws.Range("M19").Value = -830
$ws.Range("N19").Value = -2006.6666
$ws.Range("H24").Value = 1500
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 1666.6666
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1666.6666
$ws.Range("M24").Value = -830
$ws.Range("N24").Value = -2006.6666
$ws.Range("H31").Value = 1767.697
$ws.Range("I31").Value = 1761.2307
$ws.Range("K31").Value = 1761.2307
$ws.Range("M31").Value = -1466.2307
$ws.Range("H34").Value = 1767.697
$ws.Range("I34").Value = 1761.2307
$ws.Range("K34").Value = 1761.2307
$ws.Range("M34").Value = -1559.2307
$ws.Range("H107").Value = 1256.6666
$ws.Range("I107").Value = 1283.5714
$ws.Range("J107").Value = 880
$ws.Range("K107").Value = 1283.5714
$ws.Range("L107").Value = 880
$ws.Range("M107").Value = 636.4286
$ws.Range("N107").Value = -4720

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 647
$ws.Range("I5").Value = 464.9091
$ws.Range("J5").Value = 897.375
$ws.Range("K5").Value = 1394.7273
$ws.Range("L5").Value = 2692.125
$ws.Range("M5").Value = -1282.7273
$ws.Range("N5").Value = -2916.125
$ws.Range("H122").Value = 946
$ws.Range("I122").Value = 531.3333
$ws.Range("J122").Value = 1137.3846
$ws.Range("K122").Value = 4781.9997
$ws.Range("L122").Value = 10236.4614
$ws.Range("M122").Value = -2331.9997
$ws.Range("N122").Value = -15136.4614
$ws.Range("H131").Value = 11306.194
$ws.Range("I131").Value = 785
$ws.Range("J131").Value = 12341.065
$ws.Range("K131").Value = 2355
$ws.Range("L131").Value = 37023.195
$ws.Range("M131").Value = 2685
$ws.Range("N131").Value = -47103.195
$ws.Range("H132").Value = 999.8
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H135").Value = 647
$ws.Range("I135").Value = 464.9091
$ws.Range("J135").Value = 897.375
$ws.Range("K135").Value = 4184.1819
$ws.Range("L135").Value = 8076.375
$ws.Range("M135").Value = -1649.1819
$ws.Range("N135").Value = -13146.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4870.7144
$ws.Range("I70").Value = 5798.75
$ws.Range("J70").Value = 3633.3333
$ws.Range("K70").Value = 5798.75
$ws.Range("L70").Value = 3633.3333
$ws.Range("M70").Value = -5528.75
$ws.Range("N70").Value = -4173.3333
$ws.Range("H73").Value = 4870.7144
$ws.Range("I73").Value = 5798.75
$ws.Range("J73").Value = 3633.3333
$ws.Range("K73").Value = 5798.75
$ws.Range("L73").Value = 3633.3333
$ws.Range("M73").Value = -4862.75
$ws.Range("N73").Value = -5505.3333
$ws.Range("H132").Value = 875629.2
$ws.Range("I132").Value = 1132346
$ws.Range("J132").Value = 2792.3
$ws.Range("K132").Value = 3397038
$ws.Range("L132").Value = 8376.900000000001
$ws.Range("M132").Value = -3394508
$ws.Range("N132").Value = -13436.9
$ws.Range("H140").Value = 39682.715
$ws.Range("J140").Value = 39682.715
$ws.Range("L140").Value = 39682.715
$ws.Range("N140").Value = -50042.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3200.2
$ws.Range("I22").Value = 6000
$ws.Range("J22").Value = 2500.25
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 2500.25
$ws.Range("M22").Value = -5705
$ws.Range("N22").Value = -3090.25
$ws.Range("H27").Value = 3200.2
$ws.Range("I27").Value = 6000
$ws.Range("J27").Value = 2500.25
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 2500.25
$ws.Range("M27").Value = -5893
$ws.Range("N27").Value = -2714.25
$ws.Range("H40").Value = 6255.7666
$ws.Range("I40").Value = 6606.4287
$ws.Range("J40").Value = 5437.5557
$ws.Range("K40").Value = 6606.4287
$ws.Range("L40").Value = 5437.5557
$ws.Range("M40").Value = -6470.4287
$ws.Range("N40").Value = -5709.5557
$ws.Range("H93").Value = 935.5294
$ws.Range("I93").Value = 750.3077
$ws.Range("J93").Value = 1537.5
$ws.Range("K93").Value = 750.3077
$ws.Range("L93").Value = 1537.5
$ws.Range("M93").Value = 497.6923
$ws.Range("N93").Value = -4033.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2599.8
$ws.Range("I81").Value = 2499.75
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 4999.5
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -3938.5
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2599.8
$ws.Range("I84").Value = 2499.75
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 24997.5
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -19693.5
$ws.Range("N84").Value = -40608
